$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: municipio-nombre (D) and aragon (F) columns switch their
# dimension annotation from measure-specific strings to the shared
# sdmx-dimension:refArea value (matching provincia-nombre / comarca-nombre).
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "sdmx-dimension:refArea"

# Row 3: municipio-nombre (D) switches from "medida" to "dim".
$ws.Range("D3").Value = "dim"

# Row 4: municipio-nombre (D) gets a new URI column "URI-Municipio"
# (previously it held the measure type "xsd:int").
$ws.Range("D4").Value = "URI-Municipio"

# Row 4: aragon (F) switches from "skos:Concept" to a new
# "URI-Comunidad" mapping reference (no longer uses its own mapping file).
$ws.Range("F4").Value = "URI-Comunidad"

# Row 5: remove the now-unused mapping-aragon.xlsx reference in F5.
$ws.Range("F5").Clear()
